$d = $word.ActiveDocument

# 1) "LGU " + "Ginatilan" (two separate runs, split by a spell-check
#    proofErr pair) -> a single run containing "LGU Ginatilan".
#    Find/Replace across the run boundary merges the runs into one and
#    drops the proofErr markers that sit strictly inside the matched span.
$d.Content.Find.Execute("LGU Ginatilan", $true, $false, $false, $false, `
    $false, $true, 1, $false, "LGU Ginatilan", 2) | Out-Null

# 2) Remove the "IT Support Intern" run entirely, leaving the paragraph
#    (with its pPr) empty.
$r = $d.Content.Duplicate
$r.Find.Execute("IT Support Intern", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
if ($r.Find.Found) {
    $r.Delete()
}

# 3) Drop the stray "_GoBack" bookmark (both bookmarkStart/bookmarkEnd).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
